$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 2.7
$ws.Range("K2").Value = 2.1
$ws.Range("S2").Value = 1.38
$ws.Range("T2").Value = 2.82
$ws.Range("Z2").Value = 23
$ws.Range("AI2").Value = 18
$ws.Range("AO2").Value = 11
$ws.Range("AR2").Value = 60
$ws.Range("AT2").Value = 2.82
$ws.Range("AU2").Value = 6.5
$ws.Range("AX2").Value = 18.5
$ws.Range("AY2").Value = 23
$ws.Range("BA2").Value = 120
